$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Rename the three sheets (tab order stays the same, only the names change)
# ---------------------------------------------------------------------------
$wsRankable = $wb.Worksheets.Item("Sheet1")
$wsRankable.Name = "Just Rankable"

$wsExperiment = $wb.Worksheets.Item("Sheet3")
$wsExperiment.Name = "Experiment Results"

$wsResults = $wb.Worksheets.Item("Results")
$wsResults.Name = "Ranking Results"

# ---------------------------------------------------------------------------
# 2) Fix up the two existing charts (they live on "Just Rankable") so their
#    series formulas point at the renamed sheet.
# ---------------------------------------------------------------------------
$chart1 = $wsRankable.ChartObjects().Item(1).Chart
$s = $chart1.SeriesCollection().Item(1)
$s.Formula = "=SERIES('Just Rankable'!`$F`$1,'Just Rankable'!`$A`$2:`$A`$22,'Just Rankable'!`$F`$2:`$F`$22,1)"
$s = $chart1.SeriesCollection().Item(2)
$s.Formula = "=SERIES('Just Rankable'!`$G`$1,'Just Rankable'!`$A`$2:`$A`$22,'Just Rankable'!`$G`$2:`$G`$22,2)"

$chart2 = $wsRankable.ChartObjects().Item(2).Chart
$s = $chart2.SeriesCollection().Item(1)
$s.Formula = "=SERIES('Just Rankable'!`$D`$1,'Just Rankable'!`$A`$2:`$A`$22,'Just Rankable'!`$D`$2:`$D`$22,1)"
$s = $chart2.SeriesCollection().Item(2)
$s.Formula = "=SERIES('Just Rankable'!`$E`$1,'Just Rankable'!`$A`$2:`$A`$22,'Just Rankable'!`$E`$2:`$E`$22,2)"

# ---------------------------------------------------------------------------
# 3) Rebuild "Ranking Results" with the new ranking-analysis tables. The
#    sheet used to hold a small summary block; it is replaced wholesale with
#    several small crosstabs plus a per-algorithm rank-count matrix.
# ---------------------------------------------------------------------------

# Grab a cell that already carries the workbook's "s=1" style (10pt Lucida
# Console, vertically centered) so the copies below reuse that style index
# instead of minting a duplicate one.
$styleSource = $wsExperiment.Range("J2")
$styleSource.Copy() | Out-Null
$wsResults.Range("A1:Q29").PasteSpecial(-4122) | Out-Null
$wsResults.Cells.Clear()

# -- Avg rank table (A1:B8) --------------------------------------------------
$wsResults.Range("A1").Value = "AlgorithmId"
$wsResults.Range("B1").Value = "avgRank"
$avgRank = @(
    @(5, 2.1764700000000001),
    @(7, 3.1176469999999998),
    @(6, 4),
    @(4, 4),
    @(2, 4.5882350000000001),
    @(3, 5),
    @(1, 5.1176469999999998)
)
for ($i = 0; $i -lt $avgRank.Count; $i++) {
    $r = 2 + $i
    $wsResults.Range("A$r").Value = $avgRank[$i][0]
    $wsResults.Range("B$r").Value = $avgRank[$i][1]
}

# -- 1st Ranked (E1:F5) -------------------------------------------------------
$wsResults.Range("E1").Value = "AlgorithmId"
$wsResults.Range("F1").Value = "1st Ranked"
$firstRanked = @(@(5,11), @(1,4), @(2,1), @(7,1))
for ($i = 0; $i -lt $firstRanked.Count; $i++) {
    $r = 2 + $i
    $wsResults.Range("E$r").Value = $firstRanked[$i][0]
    $wsResults.Range("F$r").Value = $firstRanked[$i][1]
}

# -- 2nd Ranked (H1:I6) -------------------------------------------------------
$wsResults.Range("H1").Value = "AlgorithmId"
$wsResults.Range("I1").Value = "2nd Ranked"
$secondRanked = @(@(7,11), @(2,2), @(5,2), @(6,1), @(1,1))
for ($i = 0; $i -lt $secondRanked.Count; $i++) {
    $r = 2 + $i
    $wsResults.Range("H$r").Value = $secondRanked[$i][0]
    $wsResults.Range("I$r").Value = $secondRanked[$i][1]
}

# -- 3rd Ranking (K1:L6) -------------------------------------------------------
$wsResults.Range("K1").Value = "AlgorithmId"
$wsResults.Range("L1").Value = "3rd Ranking"
$thirdRanked = @(@(4,8), @(6,5), @(3,2), @(2,1), @(7,1))
for ($i = 0; $i -lt $thirdRanked.Count; $i++) {
    $r = 2 + $i
    $wsResults.Range("K$r").Value = $thirdRanked[$i][0]
    $wsResults.Range("L$r").Value = $thirdRanked[$i][1]
}

# -- 4th Ranked (E8:F13) ------------------------------------------------------
$wsResults.Range("E8").Value = "AlgorithmId"
$wsResults.Range("F8").Value = "4th Ranked"
$fourthRanked = @(@(6,8), @(4,4), @(2,2), @(3,2), @(1,1))
for ($i = 0; $i -lt $fourthRanked.Count; $i++) {
    $r = 9 + $i
    $wsResults.Range("E$r").Value = $fourthRanked[$i][0]
    $wsResults.Range("F$r").Value = $fourthRanked[$i][1]
}

# -- 5th Ranked (H8:I12) -------------------------------------------------------
$wsResults.Range("H8").Value = "AlgorithmId"
$wsResults.Range("I8").Value = "5th Ranked"
$fifthRanked = @(@(3,7), @(2,5), @(4,3))
for ($i = 0; $i -lt $fifthRanked.Count; $i++) {
    $r = 9 + $i
    $wsResults.Range("H$r").Value = $fifthRanked[$i][0]
    $wsResults.Range("I$r").Value = $fifthRanked[$i][1]
}
$wsResults.Range("H12").Value = 5
$wsResults.Range("I12").Value = 2

# -- 6th Ranked (K8:L14) --------------------------------------------------------
$wsResults.Range("K8").Value = "AlgorithmId"
$wsResults.Range("L8").Value = "6th Ranked"
$sixthRanked = @(@(3,6), @(2,5), @(5,2), @(6,2), @(7,1), @(4,1))
for ($i = 0; $i -lt $sixthRanked.Count; $i++) {
    $r = 9 + $i
    $wsResults.Range("K$r").Value = $sixthRanked[$i][0]
    $wsResults.Range("L$r").Value = $sixthRanked[$i][1]
}

# -- 7th Ranked (E16:F21) -------------------------------------------------------
$wsResults.Range("E16").Value = "AlgorithmId"
$wsResults.Range("F16").Value = "7th Ranked"
$seventhRanked = @(@(1,11), @(7,3), @(2,1), @(4,1), @(6,1))
for ($i = 0; $i -lt $seventhRanked.Count; $i++) {
    $r = 17 + $i
    $wsResults.Range("E$r").Value = $seventhRanked[$i][0]
    $wsResults.Range("F$r").Value = $seventhRanked[$i][1]
}

# -- Per-algorithm rank-count matrix (I16:Q23) --------------------------------
$wsResults.Range("I16").Value = "Algorithm"
$headerCols = @("J","K","L","M","N","O","P")
for ($c = 0; $c -lt 7; $c++) {
    $wsResults.Range("$($headerCols[$c])16").Value = $c + 1
}
$wsResults.Range("Q16").Value = "AVG"

$matrix = @(
    @(1, 4, 1, 0, 1, 0, 0, 11, 5.1176469999999998),
    @(2, 1, 2, 1, 2, 5, 5, 1, 4.5882350000000001),
    @(3, 0, 0, 2, 2, 7, 6, 0, 5),
    @(4, 0, 0, 8, 4, 3, 1, 1, 4),
    @(5, 11, 2, 0, 0, 2, 2, 0, 2.1764700000000001),
    @(6, 0, 1, 5, 8, 0, 2, 1, 4),
    @(7, 1, 11, 1, 0, 0, 1, 3, 3.1176469999999998)
)
for ($i = 0; $i -lt $matrix.Count; $i++) {
    $r = 17 + $i
    $row = $matrix[$i]
    $wsResults.Range("I$r").Value = $row[0]
    $wsResults.Range("J$r").Value = $row[1]
    $wsResults.Range("K$r").Value = $row[2]
    $wsResults.Range("L$r").Value = $row[3]
    $wsResults.Range("M$r").Value = $row[4]
    $wsResults.Range("N$r").Value = $row[5]
    $wsResults.Range("O$r").Value = $row[6]
    $wsResults.Range("P$r").Value = $row[7]
    $wsResults.Range("Q$r").Value = $row[8]
}

# -- cntOf1stRanked summary (C25:D29) ----------------------------------------
$wsResults.Range("C25").Value = "AlgorithmId"
$wsResults.Range("D25").Value = "cntOf1stRanked"
$cntFirst = @(@(7,10), @(1,4), @(2,1), @(5,1))
for ($i = 0; $i -lt $cntFirst.Count; $i++) {
    $r = 26 + $i
    $wsResults.Range("C$r").Value = $cntFirst[$i][0]
    $wsResults.Range("D$r").Value = $cntFirst[$i][1]
}

# ---------------------------------------------------------------------------
# 4) Formatting: reuse the workbook's "s=1" style (10pt Lucida Console,
#    vertically centred) for the small spacer cells, and apply centred /
#    centred+2dp formatting to the rank-count matrix.
# ---------------------------------------------------------------------------
$styleSource.Copy() | Out-Null
$wsResults.Range("H12").PasteSpecial(-4123) | Out-Null
$wsResults.Range("G14").PasteSpecial(-4123) | Out-Null
$wsResults.Range("F15:G15").PasteSpecial(-4123) | Out-Null
$wsResults.Range("H16").PasteSpecial(-4123) | Out-Null

$matrixHeader = $wsResults.Range("I16:Q16")
$matrixHeader.HorizontalAlignment = -4108

$matrixBody = $wsResults.Range("I17:P23")
$matrixBody.HorizontalAlignment = -4108

$avgCol = $wsResults.Range("Q17:Q23")
$avgCol.HorizontalAlignment = -4108
$avgCol.NumberFormat = "0.00"

# ---------------------------------------------------------------------------
# 5) Column width + view/selection details for "Ranking Results"
# ---------------------------------------------------------------------------
$wsResults.Columns.Item(17).ColumnWidth = 8.75

$wb.Windows.Item(1).DisplayWorkbookTabs = $true
$wsResults.Activate()
$excel.ActiveWindow.ScrollRow = 7
$wsResults.Range("C25:D29").Select() | Out-Null

Write-Output "done"
